# Update TPM-derived values on Sheet1 for rows 2-7, columns M-T (13-20)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 12.800308
$ws.Range("N2").Value = 38.400924
$ws.Range("O2").Value = 0.2552188303259509
$ws.Range("P2").Value = 0.2552188303259509
$ws.Range("Q2").Value = 0.7233240712333334
$ws.Range("R2").Value = 6.509916641100001
$ws.Range("S2").Value = 0.2552188303259509
$ws.Range("T2").Value = 0.2552188303259509

$ws.Range("M3").Value = 3.628896
$ws.Range("O3").Value = 0.07235471150338896
$ws.Range("P3").Value = 0.07235471150338896
$ws.Range("S3").Value = 0.07235471150338896
$ws.Range("T3").Value = 0.07235471150338896

$ws.Range("M4").Value = 5.520799
$ws.Range("N4").Value = 16.562397
$ws.Range("O4").Value = 0.1100764031025409
$ws.Range("P4").Value = 0.1100764031025409
$ws.Range("Q4").Value = 0.3119711501583334
$ws.Range("R4").Value = 2.807740351425
$ws.Range("S4").Value = 0.1100764031025409
$ws.Range("T4").Value = 0.1100764031025409

$ws.Range("M5").Value = 2.905047666666667
$ws.Range("N5").Value = 8.715143000000001
$ws.Range("O5").Value = 0.0579222677710411
$ws.Range("P5").Value = 0.0579222677710411
$ws.Range("Q5").Value = 0.1641594018972223
$ws.Range("R5").Value = 1.477434617075
$ws.Range("S5").Value = 0.0579222677710411
$ws.Range("T5").Value = 0.0579222677710411

$ws.Range("M6").Value = 20.19342933333333
$ws.Range("N6").Value = 60.580288
$ws.Range("O6").Value = 0.4026265160746975
$ws.Range("P6").Value = 0.4026265160746975
$ws.Range("Q6").Value = 1.141097035911111
$ws.Range("R6").Value = 10.2698733232
$ws.Range("S6").Value = 0.4026265160746975
$ws.Range("T6").Value = 0.4026265160746975

$ws.Range("M7").Value = 5.105766
$ws.Range("N7").Value = 15.317298
$ws.Range("O7").Value = 0.1018012712223807
$ws.Range("P7").Value = 0.1018012712223806
$ws.Range("Q7").Value = 0.28851832705
$ws.Range("R7").Value = 2.59666494345
$ws.Range("S7").Value = 0.1018012712223807
$ws.Range("T7").Value = 0.1018012712223806
